$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing row 3 (the reservation row) down two rows to row 5,
# opening up rows 3 and 4 for two newly "freed" slots.
$ws.Range("A3:D3").Insert()
$ws.Range("A3:D3").Insert()

# Row 2 slot id changes from 3 -> 2
$ws.Range("A2").Value = 2

# New row 3 slot (previously didn't exist)
$ws.Range("A3").Value = 3
$ws.Range("C3").Value = "Libre"
$ws.Range("D3").Value = "'4"
$ws.Range("D3").ClearFormats()

# New row 4 slot (previously didn't exist)
$ws.Range("A4").Value = 5
$ws.Range("C4").Value = "Libre"
$ws.Range("D4").Value = "'4"
$ws.Range("D4").ClearFormats()

# Row 5 (the shifted-down reservation row) slot id changes from 2 -> 1
$ws.Range("A5").Value = 1

# Update the active selection left by the user's last action
# (Excel reports the anchor cell of the drag, D4, as the active cell while
# the whole dragged range A2:D4 is highlighted.)
$ws.Range("A2:D4").Select()
